# Daily attendance processing - 2025-10-15 14:23:35
#
# Normalizes the "Recorded By" (column G) entries on the active sheet:
# whenever a cell lists "System" among the recorder names (comma
# separated), the whole list is reversed so that "System" ends up last
# instead of first (e.g. "System, dnasr281@gmail.com" becomes
# "dnasr281@gmail.com, System").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        $count = $parts.Count

        if ($count -gt 1) {
            $hasSystem = $false
            foreach ($p in $parts) {
                if ($p.ToLower() -eq "system") {
                    $hasSystem = $true
                }
            }

            if ($hasSystem) {
                $reversed = @()
                for ($i = $count - 1; $i -ge 0; $i--) {
                    $reversed += $parts[$i]
                }
                $newVal = [string]::Join(", ", $reversed)

                if ($newVal -ne $val) {
                    $cell.Value = $newVal
                }
            }
        }
    }
}
